$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''58.925.72'
$ws.Range("E2").Value = '''  -2.96%  '
$ws.Range("D3").Value = '''2.662.01'
$ws.Range("E3").Value = '''  -1.40%  '
$ws.Range("E4").Value = '''  -0.01%  '
$ws.Range("D5").Value = '''523.13'
$ws.Range("E5").Value = '''  -0.21%  '
$ws.Range("D6").Value = '''144.30'
$ws.Range("E6").Value = '''  -2.09%  '
$ws.Range("E7").Value = '''  +0.30%  '
$ws.Range("E8").Value = '''  -1.53%  '
$ws.Range("E9").Value = '''  +7.01%  '
$ws.Range("E10").Value = '''  -3.45%  '
$ws.Range("E11").Value = '''  -2.30%  '
$ws.Range("E12").Value = '''  +1.24%  '
$ws.Range("D13").Value = '''3.130.57'
$ws.Range("E13").Value = '''  -1.21%  '
$ws.Range("D14").Value = '''58.927.90'
$ws.Range("E14").Value = '''  -2.93%  '
$ws.Range("D15").Value = '''21.01'
$ws.Range("D16").Value = '''2.680.85'
$ws.Range("E16").Value = '''  -5.51%  '
$ws.Range("E17").Value = '''  -2.41%  '
$ws.Range("D18").Value = '''339.22'
$ws.Range("E18").Value = '''  -3.86%  '
$ws.Range("D19").Value = '''4.39'
$ws.Range("E19").Value = '''  -3.67%  '
$ws.Range("D20").Value = '''10.36'
$ws.Range("E20").Value = '''  -2.75%  '
$ws.Range("D21").Value = '''6.38'
$ws.Range("E21").Value = '''  +0.17%  '
$ws.Range("E22").Value = '''  -0.05%  '
$ws.Range("D23").Value = '''64.37'
$ws.Range("E23").Value = '''  +2.04%  '
$ws.Range("E24").Value = '''  -1.30%  '
$ws.Range("E25").Value = '''  -1.40%  '
$ws.Range("D26").Value = '''0.997'
$ws.Range("E26").Value = '''  +0.56%  '
$ws.Range("E27").Value = '''  -2.62%  '
$ws.Range("D28").Value = '''7.14'
$ws.Range("E28").Value = '''  -2.18%  '
$ws.Range("D29").Value = '''6.66'
$ws.Range("E29").Value = '''  -2.78%  '
$ws.Range("E30").Value = '''  +0.08%  '
$ws.Range("E31").Value = '''  -0.32%  '
$ws.Range("D32").Value = '''18.89'
$ws.Range("E32").Value = '''  -1.47%  '
$ws.Range("D33").Value = '''150.53'
$ws.Range("E33").Value = '''  +1.72%  '
$ws.Range("D34").Value = '''4.15'
$ws.Range("E34").Value = '''  -3.75%  '
$ws.Range("E35").Value = '''  -5.66%  '
$ws.Range("D36").Value = '''0.900'
$ws.Range("E36").Value = '''  -5.64%  '
$ws.Range("D37").Value = '''0.871'
$ws.Range("E37").Value = '''  -1.50%  '
$ws.Range("D38").Value = '''36.85'
$ws.Range("E38").Value = '''  -0.37%  '
$ws.Range("E39").Value = '''  -6.21%  '
$ws.Range("D40").Value = '''3.58'
$ws.Range("E40").Value = '''  -3.68%  '
$ws.Range("D41").Value = '''0.615'
$ws.Range("E41").Value = '''  -0.10%  '
$ws.Range("E42").Value = '''  +0.40%  '
$ws.Range("D43").Value = '''275.43'
$ws.Range("E43").Value = '''  -2.99%  '
$ws.Range("D44").Value = '''19.81'
$ws.Range("E44").Value = '''  -1.58%  '
$ws.Range("E45").Value = '''  -2.25%  '
$ws.Range("D46").Value = '''10.67'
$ws.Range("E46").Value = '''  +1.98%  '
$ws.Range("D47").Value = '''0.0533'
$ws.Range("E47").Value = '''  -1.42%  '
$ws.Range("D48").Value = '''2.051.36'
$ws.Range("E48").Value = '''  -4.48%  '
$ws.Range("D49").Value = '''4.73'
$ws.Range("E49").Value = '''  -4.00%  '
$ws.Range("E50").Value = '''  -3.27%  '
$ws.Range("D51").Value = '''18.78'
$ws.Range("E51").Value = '''  -3.35%  '
